$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 39670
$ws.Range("J109").Value = 39670
$ws.Range("L109").Value = 39670
$ws.Range("N109").Value = -42444
$ws.Range("H117").Value = 45659.332
$ws.Range("J117").Value = 45659.332
$ws.Range("L117").Value = 45659.332
$ws.Range("N117").Value = -54837.332
$ws.Range("H127").Value = 1387.875
$ws.Range("J127").Value = 1433.7333
$ws.Range("L127").Value = 4301.199900000001
$ws.Range("N127").Value = -14221.1999
$ws.Range("H128").Value = 42904.715
$ws.Range("J128").Value = 42904.715
$ws.Range("L128").Value = 42904.715
$ws.Range("N128").Value = -52864.715
$ws.Range("H129").Value = 1220.6703
$ws.Range("J129").Value = 1160.0476
$ws.Range("L129").Value = 3480.142800000001
$ws.Range("N129").Value = -13480.1428
$ws.Range("H131").Value = 3209.9443
$ws.Range("I131").Value = 3793.3333
$ws.Range("J131").Value = 3093.2666
$ws.Range("K131").Value = 11379.9999
$ws.Range("L131").Value = 9279.799800000001
$ws.Range("M131").Value = -6339.999899999999
$ws.Range("N131").Value = -19359.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 12000
$ws.Range("J39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -13040
$ws.Range("H80").Value = 60046
$ws.Range("J80").Value = 60046
$ws.Range("L80").Value = 60046
$ws.Range("N80").Value = -62042
$ws.Range("H83").Value = 60046
$ws.Range("J83").Value = 60046
$ws.Range("L83").Value = 180138
$ws.Range("N83").Value = -190122
$ws.Range("H117").Value = 44098.5
$ws.Range("J117").Value = 44098.5
$ws.Range("L117").Value = 44098.5
$ws.Range("N117").Value = -53276.5
$ws.Range("H120").Value = 38022.4
$ws.Range("J120").Value = 38022.4
$ws.Range("L120").Value = 38022.4
$ws.Range("N120").Value = -47698.4
$ws.Range("H122").Value = 4309.75
$ws.Range("I122").Value = 4929.6665
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 14788.9995
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -12338.9995
$ws.Range("N122").Value = -12250
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H130").Value = 34562.5
$ws.Range("J130").Value = 34562.5
$ws.Range("L130").Value = 34562.5
$ws.Range("N130").Value = -44602.5
$ws.Range("H131").Value = 48107.4
$ws.Range("J131").Value = 48107.4
$ws.Range("L131").Value = 48107.4
$ws.Range("N131").Value = -58187.4
$ws.Range("H132").Value = 13890607
$ws.Range("I132").Value = 21740324
$ws.Range("J132").Value = 2646.8462
$ws.Range("K132").Value = 65220972
$ws.Range("L132").Value = 7940.5386
$ws.Range("M132").Value = -65218442
$ws.Range("N132").Value = -13000.5386
$ws.Range("H133").Value = 37720.848
$ws.Range("J133").Value = 37720.848
$ws.Range("L133").Value = 37720.848
$ws.Range("N133").Value = -42780.848

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 33625
$ws.Range("J38").Value = 33625
$ws.Range("L38").Value = 33625
$ws.Range("N38").Value = -34457
$ws.Range("H124").Value = 49992
$ws.Range("J124").Value = 49992
$ws.Range("L124").Value = 49992
$ws.Range("N124").Value = -59812
$ws.Range("H125").Value = 50352
$ws.Range("J125").Value = 50352
$ws.Range("L125").Value = 50352
$ws.Range("N125").Value = -60192
$ws.Range("H137").Value = 26940
$ws.Range("J137").Value = 26940
$ws.Range("L137").Value = 26940
$ws.Range("N137").Value = -37140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49883.5
$ws.Range("J20").Value = 49883.5
$ws.Range("L20").Value = 49883.5
$ws.Range("N20").Value = -50355.5
$ws.Range("H30").Value = 49883.5
$ws.Range("J30").Value = 49883.5
$ws.Range("L30").Value = 49883.5
$ws.Range("N30").Value = -50065.5
$ws.Range("H99").Value = 1840
$ws.Range("I99").Value = 1777.5
$ws.Range("J99").Value = 1865
$ws.Range("K99").Value = 1777.5
$ws.Range("L99").Value = 1865
$ws.Range("M99").Value = -279.5
$ws.Range("N99").Value = -4861
$ws.Range("H116").Value = 46744
$ws.Range("J116").Value = 46744
$ws.Range("L116").Value = 46744
$ws.Range("N116").Value = -55922
$ws.Range("H126").Value = 1840
$ws.Range("I126").Value = 1777.5
$ws.Range("J126").Value = 1865
$ws.Range("K126").Value = 5332.5
$ws.Range("L126").Value = 5595
$ws.Range("M126").Value = -2862.5
$ws.Range("N126").Value = -10535
$ws.Range("H128").Value = 49883.5
$ws.Range("J128").Value = 49883.5
$ws.Range("L128").Value = 49883.5
$ws.Range("N128").Value = -59843.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2003.1111
$ws.Range("I102").Value = 1770.3077
$ws.Range("K102").Value = 1770.3077
$ws.Range("M102").Value = -148.3077000000001
$ws.Range("H114").Value = 45904.668
$ws.Range("J114").Value = 45904.668
$ws.Range("L114").Value = 45904.668
$ws.Range("N114").Value = -54582.668
$ws.Range("H116").Value = 49734
$ws.Range("J116").Value = 49734
$ws.Range("L116").Value = 49734
$ws.Range("N116").Value = -58912
$ws.Range("H122").Value = 1135.4445
$ws.Range("I122").Value = 1164.619
$ws.Range("J122").Value = 1033.3334
$ws.Range("K122").Value = 3493.857
$ws.Range("L122").Value = 3100.0002
$ws.Range("M122").Value = -1043.857
$ws.Range("N122").Value = -8000.0002
$ws.Range("H126").Value = 6626.7827
$ws.Range("I126").Value = 11560
$ws.Range("J126").Value = 2104.6667
$ws.Range("K126").Value = 34680
$ws.Range("L126").Value = 6314.000100000001
$ws.Range("M126").Value = -32210
$ws.Range("N126").Value = -11254.0001
$ws.Range("H130").Value = 50544
$ws.Range("J130").Value = 50544
$ws.Range("L130").Value = 50544
$ws.Range("N130").Value = -60584

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2249.7568
$ws.Range("I7").Value = 1970.375
$ws.Range("K7").Value = 1970.375
$ws.Range("M7").Value = -1858.375
$ws.Range("H40").Value = 2135
$ws.Range("I40").Value = 2193
$ws.Range("J40").Value = 1700
$ws.Range("K40").Value = 2193
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = -2057
$ws.Range("N40").Value = -1972
$ws.Range("H111").Value = 43968
$ws.Range("J111").Value = 43968
$ws.Range("L111").Value = 43968
$ws.Range("N111").Value = -52148
$ws.Range("H116").Value = 50668
$ws.Range("J116").Value = 50668
$ws.Range("L116").Value = 50668
$ws.Range("N116").Value = -59846
$ws.Range("H122").Value = 1978.5714
$ws.Range("I122").Value = 1583.3334
$ws.Range("J122").Value = 2275
$ws.Range("K122").Value = 4750.0002
$ws.Range("L122").Value = 6825
$ws.Range("M122").Value = -2300.0002
$ws.Range("N122").Value = -11725
$ws.Range("H126").Value = 2249.7568
$ws.Range("I126").Value = 1970.375
$ws.Range("K126").Value = 5911.125
$ws.Range("M126").Value = -3441.125
$ws.Range("H127").Value = 25974
$ws.Range("J127").Value = 42233.2
$ws.Range("L127").Value = 42233.2
$ws.Range("N127").Value = -52153.2
$ws.Range("H130").Value = 48329
$ws.Range("J130").Value = 48329
$ws.Range("L130").Value = 48329
$ws.Range("N130").Value = -58369
$ws.Range("H135").Value = 45747.5
$ws.Range("J135").Value = 45747.5
$ws.Range("L135").Value = 45747.5
$ws.Range("N135").Value = -55887.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1505025.5
$ws.Range("I122").Value = 4082439.8
$ws.Range("K122").Value = 12247319.4
$ws.Range("M122").Value = -12244869.4
$ws.Range("H126").Value = 9804321
$ws.Range("I126").Value = 9804321
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 29412963
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -29410493
$ws.Range("N126").ClearContents()
